# Generate Report for Handoff
# b.md has been newly handed off for localization (zh-cn, de-de): update its
# status to "Ready for handoff" and record the new handoff file/datetime.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: just the two status cells for b.md change.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("D3").Value = "2016-03-09 08:02:13"

# The hyperlink in C3 needs its display text updated to point at the new
# handoff file name, while its target address stays the same. This runtime
# only supports clearing *all* hyperlinks on a sheet at once, so clear them
# and re-create every one of them (address unchanged, display text updated
# only for C3).
$zh.Range("A1").Hyperlinks.Delete()

$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d84e0381bf515289d4095f2c86b89a80ca90c299/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/22368cfc249d46b42854d502738d4488bce096b2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3d4fc35eb750f888ef84a3b7e27c39a4dfff1777/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc689d5cd57948b550d45b4b0725b72b2bf91373/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d84e0381bf515289d4095f2c86b89a80ca90c299/e2e/b.md", "", "", "b.md")
$zh.Hyperlinks.Add($zh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/22368cfc249d46b42854d502738d4488bce096b2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3d4fc35eb750f888ef84a3b7e27c39a4dfff1777/e2e/a.md", "", "", "a.md")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/dc689d5cd57948b550d45b4b0725b72b2bf91373/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d84e0381bf515289d4095f2c86b89a80ca90c299/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B3").Value = "Ready for handoff"
$de.Range("D3").Value = "2016-03-09 08:02:20"

$de.Range("A1").Hyperlinks.Delete()

$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d84e0381bf515289d4095f2c86b89a80ca90c299/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e8d692954314a9c3b9ab633941c78d5ef009fa1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/12ac3220f4305b18384e1f1cf6e3e7f486132cfd/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63ea7edfc43a1bc82a09e2b552951ec947c4da9d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d84e0381bf515289d4095f2c86b89a80ca90c299/e2e/b.md", "", "", "b.md")
$de.Hyperlinks.Add($de.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e8d692954314a9c3b9ab633941c78d5ef009fa1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/12ac3220f4305b18384e1f1cf6e3e7f486132cfd/e2e/a.md", "", "", "a.md")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63ea7edfc43a1bc82a09e2b552951ec947c4da9d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d84e0381bf515289d4095f2c86b89a80ca90c299/.localization-config", "", "", ".localization-config")

$wb.Save()
